# feat: add 2022-Q1 data
#
# The workbook has quarterly "fund holdings" sheets (2020-Q4 .. 2021-Q4) plus
# a "总计" (totals) summary sheet. This change:
#   1. Turns the current "总计" sheet into the new "2022-Q1" holdings sheet
#      (same column layout as the other quarterly sheets: 基金代码, 基金名称,
#      基金规模, 股票总仓位, 仓位占比, 持有市值(亿元), 仓位排名).
#   2. Appends a brand new "总计" sheet at the end of the workbook with the
#      same totals layout as before, plus a new first data row for 2022-Q1.

$wb = $excel.ActiveWorkbook

$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# Step 1: convert the existing "总计" sheet into the "2022-Q1" sheet.
# ---------------------------------------------------------------------
$totalWs = $wb.Worksheets.Item("总计")

# Use another quarterly sheet as a style template for the header row and
# the bold/bordered "index" column (style index 2 in the original file).
$template = $wb.Worksheets.Item("2021-Q4")

# Wipe everything the old totals sheet had (values + formatting) so no
# stray formatting/content from the 4-column totals table survives.
$totalWs.Cells.Clear()

# Copy the whole used block of the template sheet (header + data rows) onto
# the target so column widths/row styles/borders match the other quarterly
# sheets, then we will overwrite the actual cell values below.
$template.Range("A1:H23").Copy($totalWs.Range("A1"))

# Clear any leftover rows below the 5 rows (1 header + 4 data) we need.
$totalWs.Range("A6:H23").Clear()

$totalWs.Name = "2022-Q1"

# Header row.
$totalWs.Cells.Item(1,2).Value = "基金代码"
$totalWs.Cells.Item(1,3).Value = "基金名称"
$totalWs.Cells.Item(1,4).Value = "基金规模"
$totalWs.Cells.Item(1,5).Value = "股票总仓位"
$totalWs.Cells.Item(1,6).Value = "仓位占比"
$totalWs.Cells.Item(1,7).Value = "持有市值(亿元)"
$totalWs.Cells.Item(1,8).Value = "仓位排名"

# Fund holdings data rows (column A is the existing 0-based index column,
# already populated 0..21 by the template copy; rows 2-5 are the ones we
# need, so nothing to change there).
$fundRows = @(
    @(0, "006371", "长安鑫盈灵活配置混合A", "10.52", "88.39", "3.09", "0.3251", 10),
    @(1, "006323", "合煦智远嘉选混合A",     "1.54",  "79.45", "6.01", "0.0926", 1),
    @(2, "006372", "长安鑫盈灵活配置混合C", "1.60",  "88.39", "3.09", "0.0494", 10),
    @(3, "006324", "合煦智远嘉选混合C",     "0.59",  "79.45", "6.01", "0.0355", 1)
)

# Force the numeric-looking text columns (B, D, E, F, G) to be written as
# text rather than being auto-parsed into numbers (which would drop leading
# zeros in fund codes and turn D/E/F/G into numeric cells).
for ($i = 0; $i -lt $fundRows.Length; $i++) {
    $r = $i + 2
    $totalWs.Range("B$r" + ":G$r").NumberFormat = "@"
}

for ($i = 0; $i -lt $fundRows.Length; $i++) {
    $r = $i + 2
    $row = $fundRows[$i]
    $totalWs.Cells.Item($r,1).Value = $row[0]
    $totalWs.Cells.Item($r,2).Value = "'" + $row[1]
    $totalWs.Cells.Item($r,3).Value = $row[2]
    $totalWs.Cells.Item($r,4).Value = "'" + $row[3]
    $totalWs.Cells.Item($r,5).Value = "'" + $row[4]
    $totalWs.Cells.Item($r,6).Value = "'" + $row[5]
    $totalWs.Cells.Item($r,7).Value = "'" + $row[6]
    $totalWs.Cells.Item($r,8).Value = $row[7]
}

# The apostrophe-prefix trick above stamps a "quote prefix" flag onto the
# cell's style, which would leave the data cells on a different style index
# than the rest of the workbook (they should stay on the plain, no-style
# format like every other data cell). Re-paste the plain formatting from a
# genuine no-style cell over the values we just wrote to strip that back off
# without touching the text we just stored.
$plainTemplate = $template.Range("B2:H5")
$plainTemplate.Copy()
$totalWs.Range("B2:H5").PasteSpecial($xlPasteFormats)

$totalWs.Range("A1").Select()

# ---------------------------------------------------------------------
# Step 2: append a brand new "总计" sheet at the end of the workbook.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newTotalWs = $wb.Worksheets.Add($null, $lastSheet)
$newTotalWs.Name = "总计"

# Re-use the old totals sheet's layout/style by copying it from the
# "2022-Q1" sheet's current neighbour-in-spirit: the 2021-Q4 template has
# the right header style (bold + border) for column B:D / A, so borrow the
# same style cells used on the original totals sheet (still index 2).
$template.Range("A1:D6").Copy($newTotalWs.Range("A1"))
$newTotalWs.Range("A8:D23").Clear()

# Header row for the totals sheet.
$newTotalWs.Cells.Item(1,2).Value = "日期"
$newTotalWs.Cells.Item(1,3).Value = "持有数量(只)"
$newTotalWs.Cells.Item(1,4).Value = "持有市值(亿元)"
# Column A has no header.
$newTotalWs.Cells.Item(1,1).ClearContents()

$totalRows = @(
    @(0, "2022-Q1", 4,  0.5),
    @(1, "2021-Q4", 22, 3.42),
    @(2, "2021-Q3", 35, 5.89),
    @(3, "2021-Q2", 20, 2.33),
    @(4, "2021-Q1", 2,  0.01),
    @(5, "2020-Q4", 1,  0)
)

foreach ($row in $totalRows) {
    $r = [int]$row[0] + 2
    $newTotalWs.Cells.Item($r,1).Value = $row[0]
    $newTotalWs.Cells.Item($r,2).Value = $row[1]
    $newTotalWs.Cells.Item($r,3).Value = $row[2]
    $newTotalWs.Cells.Item($r,4).Value = $row[3]
}

# Make sure there's nothing left over past row 7 (the template range we
# copied above covered rows 1-6 only, but be defensive).
$newTotalWs.Range("A8:D50").Clear()

$newTotalWs.Range("A1").Select()

Write-Output "2022-Q1 sheet added; 总计 sheet rebuilt"
